$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet values (docs/CodeSystem-duo-codes.xlsx) ---
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://purl.obolibrary.org/obo/duo.owl"
# Status
$meta.Range("B6").Value = "draft"
# Date
$meta.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Restore/apply wrap-text + top-vertical alignment on every used cell   ---
# --- so the cell formats carry an explicit alignment flag (applyAlignment) ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $used.WrapText = $true
    $used.VerticalAlignment = -4160
}
